$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Hoja1" - the workbook's only sheet

# 1) Header cell K2 ("Ganancia  por plato ") is retitled to
#    "Ganancia por Unidad de plato " (column now holds the per-unit gain,
#    while M2 keeps the existing "Ganancia por plato Vendido" label).
$ws.Range("K2").Value = "Ganancia por Unidad de plato "

# 2) Row 3 (the merged blank header row under the titles) grows taller
#    to fit the now-longer wrapped header text.
$ws.Rows.Item(3).RowHeight = 33.75

# 3) Update the saved selection/active cell shown when the sheet is reopened.
$ws.Range("L9").Select() | Out-Null
